# Weekly update: insert a new price-report row for "Haba" at the top of
# the data table (row 71), pushing the existing rows down by one.
# The last existing row (old row 82) is duplicated into the new row 83,
# matching the source data's rollover behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 71; this shifts rows 71:82 down to 72:83
# and copies formatting (including the date cell style) from row 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with this week's figures.
$ws.Range("A71").Value = 7
$ws.Range("B71").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C71").Value = 'Ñuble'
$ws.Range("D71").Value = 45142
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = 100112026
$ws.Range("G71").Value = 'Haba'
$ws.Range("H71").Value = 'Sin especificar'
$ws.Range("I71").Value = 'Primera'
$ws.Range("J71").Value = 30
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = 15000
$ws.Range("N71").Value = '$/saco 25 kilos'
$ws.Range("O71").Value = 'Provincia de Diguillín'
$ws.Range("P71").Value = 600
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = 'Hortaliza'
